$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: split the "insight" run so that "is" is wrapped in proofErr
# gramStart/gramEnd marks (Word's grammar checker flagging "there is 5
# fingers").
# ---------------------------------------------------------------------------
$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*problem listed is to figure out*") {
        $target1 = $p
        break
    }
}

if ($target1 -eq $null) {
    throw "Could not locate paragraph 1 (problem listed)"
}

$rng1 = $target1.Range
$xml1 = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p w14:paraId="38831D4A" w14:textId="5375E36E" w:rsidR="00A357D0" w:rsidRDefault="00A357D0" w:rsidP="00F85C62" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:r><w:tab/><w:t xml:space="preserve">The problem listed is to figure out which fingers she would stop at given the counting system told in the problem.  </w:t></w:r><w:r w:rsidR="00DC0713"><w:t xml:space="preserve">The insight that I noticed while reading the problem is that the fingers would keep a consistent pattern of their numbering because there </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>is</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 5 fingers in a 10 count system. </w:t></w:r><w:r w:rsidR="00DA3BF8"><w:t xml:space="preserve">The overall goal is to figure out which finger she would stop on at the various intervals given. </w:t></w:r></w:p></w:body>
'@
$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: append a new closing paragraph (separated by a blank paragraph)
# at the very end of the document, moving the _GoBack bookmark onto it.
# ---------------------------------------------------------------------------
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*solution I would use is examining the pattern*") {
        $target2 = $p
    }
}

if ($target2 -eq $null) {
    throw "Could not locate paragraph 2 (solution I would use)"
}

$rng2 = $target2.Range
$xml2 = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p w14:paraId="372EAE88" w14:textId="5DBC19E2" w:rsidR="00E97805" w:rsidRPr="00A357D0" w:rsidRDefault="00E97805" w:rsidP="00F85C62" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:r><w:tab/><w:t xml:space="preserve">The solution I would use is examining the pattern. Because the pattern is 5 fingers in a 10 count the pattern remains the same, regardless of how many times it is repeated. This would work in all cases of this problem, provided the amount of fingers counted and the count itself remained the same. </w:t></w:r></w:p></w:body>
'@
$rng2.InsertXML($xml2)

# Re-resolve the (now bookmark-free) paragraph and append the two new
# paragraphs after it.
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*solution I would use is examining the pattern*") {
        $target2 = $p
    }
}
$rng2 = $target2.Range
$rng2.InsertParagraphAfter()

$blankIndex = $target2.Index + 1
$blankPara = $d.Paragraphs($blankIndex)
$xmlBlank = '<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p/></w:body>'
$blankPara.Range.InsertXML($xmlBlank)
$blankPara = $d.Paragraphs($blankIndex)
$blankPara.Range.InsertParagraphAfter()

$finalIndex = $blankIndex + 1
$finalPara = $d.Paragraphs($finalIndex)
$rngFinal = $finalPara.Range
$xmlFinal = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:tab/><w:t xml:space="preserve">My final solution is to simply infer results from the pattern. No matter how many times the process is repeated the fingers that the count lands on would remain the same. I attempted to overthink this problem as well, thinking about math equations and such. In the end, the obvious solution became the simplest. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body>
'@
$rngFinal.InsertXML($xmlFinal)

Write-Host "Edit complete. Total paragraphs: $($d.Paragraphs.Count)"
